$d = $word.ActiveDocument

function Trimmed($range) {
    return $range.Text.TrimEnd([char]13, [char]10)
}

# Build a 1-based index -> paragraph text map (trimmed of the trailing
# paragraph mark) once, so we can locate the paragraphs we need by their
# content instead of relying on brittle fixed positions.
$paraList = @()
foreach ($para in $d.Paragraphs) {
    $paraList += $para
}

# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
for ($i = 0; $i -lt $paraList.Count; $i++) {
    $txt = Trimmed $paraList[$i].Range
    if ($txt -eq "September 19, 2025") {
        $paraList[$i].Range.Text = "September 21, 2025"
        break
    }
}

# 2. Split the recipient mailing address paragraph
#    "221 Caldwell Avenue, Los Gatos CA 95032" into two paragraphs:
#      "221 Caldwell Avenue"
#      "Los Gatos, CA 95032"
for ($i = 0; $i -lt $paraList.Count; $i++) {
    $addrPara = $paraList[$i]
    if ($addrPara.Range.Information(12)) {
        # Skip any occurrence that lives inside a table (e.g. the
        # "PROPERTY ADDRESS:" row) - only the standalone mailing-address
        # paragraph should be split.
        continue
    }
    $txt = Trimmed $addrPara.Range
    if ($txt -eq "221 Caldwell Avenue, Los Gatos CA 95032") {
        $addrPara.Range.InsertParagraphAfter()
        # Re-fetch paragraphs since the collection/indices shifted.
        $line1 = $d.Paragraphs($i + 1)
        $line2 = $d.Paragraphs($i + 2)
        $line1.Range.Text = "221 Caldwell Avenue"
        $line2.Range.Text = "Los Gatos, CA 95032"
        break
    }
}

# 3. Remove the blank "No Spacing" paragraph that immediately follows the
#    "Board of Directors" signature line.
$targetIndex = -1
for ($i = 0; $i -lt $paraList.Count; $i++) {
    $txt = Trimmed $paraList[$i].Range
    if ($txt -like "*Board of Directors*") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -ge 0) {
    $blank = $d.Paragraphs($targetIndex + 2)
    if ((Trimmed $blank.Range) -eq "") {
        $blank.Range.Delete()
    }
}
